# Update crypto price/volume table to the latest scraped values
# (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the cell to stay a text value (avoids Excel re-interpreting
    # numeric-looking strings such as "64.109.88" or "1.00" as numbers),
    # then drop back to the default "Normal" style so no new cell-level
    # style index is left behind.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "64.109.88"
Set-TextCell "E2" "  -0.09%  "
# Row 3
Set-TextCell "D3" "2.759.40"
Set-TextCell "E3" "  +0.98%  "
# Row 4
Set-TextCell "E4" "  +0.06%  "
# Row 5
Set-TextCell "D5" "576.69"
Set-TextCell "E5" "  -0.71%  "
# Row 6
Set-TextCell "D6" "159.30"
Set-TextCell "E6" "  +0.67%  "
# Row 7
Set-TextCell "E7" "  +0.19%  "
# Row 8
Set-TextCell "E8" "  -3.70%  "
# Row 9
Set-TextCell "E9" "  -1.92%  "
# Row 10
Set-TextCell "D10" "0.166"
Set-TextCell "E10" "  +4.23%  "
# Row 11
Set-TextCell "D11" "0.387"
Set-TextCell "E11" "  -1.87%  "
# Row 12
Set-TextCell "E12" "  -16.24%  "
# Row 13
Set-TextCell "D13" "3.247.59"
Set-TextCell "E13" "  +0.73%  "
# Row 14
Set-TextCell "D14" "26.95"
Set-TextCell "E14" "  -1.54%  "
# Row 15
Set-TextCell "D15" "63.707.72"
Set-TextCell "E15" "  -0.53%  "
# Row 16
Set-TextCell "D16" "0.0000152"
Set-TextCell "E16" "  -2.42%  "
# Row 17
Set-TextCell "D17" "2.761.88"
Set-TextCell "E17" "  +0.27%  "
# Row 18
Set-TextCell "E18" "  +0.30%  "
# Row 19
Set-TextCell "D19" "4.87"
Set-TextCell "E19" "  -1.83%  "
# Row 20
Set-TextCell "D20" "356.46"
Set-TextCell "E20" "  -2.31%  "
# Row 21
Set-TextCell "D21" "6.76"
Set-TextCell "E21" "  -3.89%  "
# Row 22
Set-TextCell "D22" "1.00"
Set-TextCell "E22" "  +0.71%  "
# Row 23
Set-TextCell "D23" "0.534"
Set-TextCell "E23" "  -1.56%  "
# Row 24
Set-TextCell "D24" "65.50"
Set-TextCell "E24" "  -2.24%  "
# Row 25
Set-TextCell "E25" "  -1.25%  "
# Row 26
Set-TextCell "E26" "  -0.73%  "
# Row 27
Set-TextCell "E27" "  +0.15%  "
# Row 28
Set-TextCell "D28" "0.0₃0910"
Set-TextCell "E28" "  -1.35%  "
# Row 29
Set-TextCell "B29" "Aptos"
Set-TextCell "C29" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D29" "7.29"
Set-TextCell "E29" "  +0.41%  "
# Row 30
Set-TextCell "B30" "PancakeSwap"
Set-TextCell "C30" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D30" "1.96"
Set-TextCell "E30" "  -3.12%  "
# Row 31
Set-TextCell "D31" "1.25"
Set-TextCell "E31" "  -0.53%  "
# Row 32
Set-TextCell "D32" "169.57"
Set-TextCell "E32" "  -2.40%  "
# Row 33
Set-TextCell "D33" "20.29"
Set-TextCell "E33" "  -1.95%  "
# Row 34
Set-TextCell "D34" "4.93"
Set-TextCell "E34" "  -0.67%  "
# Row 35
Set-TextCell "D35" "1.49"
Set-TextCell "E35" "  +2.04%  "
# Row 36
Set-TextCell "E36" "  +0.09%  "
# Row 37
Set-TextCell "E37" "  -1.03%  "
# Row 38
Set-TextCell "E38" "  -0.15%  "
# Row 39
Set-TextCell "B39" "Bittensor"
Set-TextCell "C39" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D39" "345.32"
Set-TextCell "E39" "  +1.72%  "
# Row 40
Set-TextCell "B40" "RenderToken"
Set-TextCell "C40" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D40" "6.32"
Set-TextCell "E40" "  +0.46%  "
# Row 41
Set-TextCell "D41" "4.21"
Set-TextCell "E41" "  -2.37%  "
# Row 42
Set-TextCell "D42" "39.24"
Set-TextCell "E42" "  -0.65%  "
# Row 43
Set-TextCell "D43" "21.44"
Set-TextCell "E43" "  -2.21%  "
# Row 44
Set-TextCell "D44" "21.82"
Set-TextCell "E44" "  -2.78%  "
# Row 45
Set-TextCell "D45" "0.0591"
Set-TextCell "E45" "  -1.99%  "
# Row 46
Set-TextCell "B46" "VeChain"
Set-TextCell "C46" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D46" "0.0255"
Set-TextCell "E46" "  -1.77%  "
# Row 47
Set-TextCell "B47" "Mantle"
Set-TextCell "C47" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D47" "0.633"
Set-TextCell "E47" "  -2.24%  "
# Row 48
Set-TextCell "E48" "  -0.83%  "
# Row 49
Set-TextCell "D49" "135.79"
Set-TextCell "E49" "  -1.64%  "
# Row 50
Set-TextCell "D50" "0.998"
Set-TextCell "E50" "  -0.06%  "
# Row 51
Set-TextCell "E51" "  +0.25%  "
